# Redeem points 71277628 766.0
#
# The sheet logs point redemptions as (phone, points, timestamp) rows.
# This edit:
#   1. Normalizes row 9's phone number (A9) from a text value to a real
#      number, matching every other phone-number cell in column A.
#   2. Appends a new redemption row (row 10) for phone 71277628 redeeming
#      766 points, timestamped 2025-08-18T16:53:54. The phone number is
#      written as text (as it originally arrived for row 9) rather than a
#      number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 9: store the phone number as a genuine number.
$ws.Range("A9").Value = 71277628

# 2) Row 10: new redemption entry.
#    Leading apostrophe forces Excel to keep the numeric-looking phone
#    number as text instead of auto-converting it to a number.
$ws.Range("A10").Value = "'71277628"
$ws.Range("B10").Value = 766
$ws.Range("C10").Value = "2025-08-18T16:53:54"
